# Add five new numbered paragraphs right after the existing "Hello"
# paragraph, while leaving the trailing _GoBack bookmark wrapping the
# last of the new paragraphs (as it originally wrapped the end of the
# "Hello" paragraph).
#
# The engine keeps a bookmark bound to whatever paragraph node
# originally hosted it, so rather than appending text after "Hello" (which
# would leave the bookmark stuck on the "Hello" paragraph), we insert the
# new lines *before* the start of the "Hello" paragraph's range. That
# pushes "Hello" down to become the final paragraph of the block, still
# carrying the bookmark, and then we simply overwrite that now-last
# paragraph's text with the final line of the list.

$d = $word.ActiveDocument

$firstPara = $d.Paragraphs.First
$r = $firstPara.Range
$r.Collapse(1)  # wdCollapseStart

$newLines = @(
    "Hello",
    "1. Help people to understand and get info from data ASAP.",
    "2. Data analysis",
    "3. Present information to business users",
    "4. Data management"
)

$r.InsertBefore(($newLines -join "`r") + "`r")

$lastPara = $d.Paragraphs.Last
$lastPara.Range.Text = "5. A connection between back end and front end"
